$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario")

$ws.Range("B2").Value = 267.289004
$ws.Range("C2").Value = 527.785518

$ws.Range("B3").Value = 236.473001
$ws.Range("C3").Value = 488.515772

$ws.Range("B4").Value = 191.882325
$ws.Range("C4").Value = 414.575915

$ws.Range("B5").Value = 150.330923
$ws.Range("C5").Value = 323.723967

$ws.Range("B6").Value = 112.816196
$ws.Range("C6").Value = 231.341876

$ws.Range("B7").Value = 88.552161
$ws.Range("C7").Value = 157.35794

$ws.Range("B8").Value = 71.030787
$ws.Range("C8").Value = 89.632421
